# Add CSV export, backups, and encrypted user storage (#2)
#
# This edit reflects new inventory-management activity recorded by the app:
# a few product stock levels were adjusted (sales processed), and a new
# user session ("Blestro") logged in and made several sales, all of which
# got appended to the "Historia" (history) sheet.

$wb = $excel.ActiveWorkbook

# --- Update stock levels on the "Inventario" sheet ---
$inv = $wb.Worksheets.Item("Inventario")
$inv.Range("C3").Value = 50
$inv.Range("C5").Value = 33
$inv.Range("C6").Value = 140
$inv.Range("C8").Value = 4
$inv.Range("C9").Value = 99862

# --- Append new movement rows to the "Historia" sheet ---
$hist = $wb.Worksheets.Item("Historia")

$hist.Range("A21").Value = "2025-08-04 20:28:40"
$hist.Range("B21").Value = "Blestro"
$hist.Range("C21").Value = "-"
$hist.Range("D21").Value = 0
$hist.Range("E21").Value = 0
$hist.Range("F21").Value = 0
$hist.Range("G21").Value = "Inicio de sesión"

$hist.Range("A22").Value = "2025-08-04 20:31:26"
$hist.Range("B22").Value = "Blestro"
$hist.Range("C22").Value = "Coca-Cola 500ml"
$hist.Range("D22").Value = -1
$hist.Range("E22").Value = 1000
$hist.Range("F22").Value = 1000
$hist.Range("G22").Value = "Venta"

$hist.Range("A23").Value = "2025-08-04 20:31:31"
$hist.Range("B23").Value = "Blestro"
$hist.Range("C23").Value = "Fernet Branca"
$hist.Range("D23").Value = -1
$hist.Range("E23").Value = 4500
$hist.Range("F23").Value = 4500
$hist.Range("G23").Value = "Venta"

$hist.Range("A24").Value = "2025-08-04 20:31:40"
$hist.Range("B24").Value = "Blestro"
$hist.Range("C24").Value = "nalga de tom"
$hist.Range("D24").Value = -22
$hist.Range("E24").Value = 115
$hist.Range("F24").Value = 2530
$hist.Range("G24").Value = "Venta"

$hist.Range("A25").Value = "2025-08-04 20:31:45"
$hist.Range("B25").Value = "Blestro"
$hist.Range("C25").Value = "Coca-Cola 500ml"
$hist.Range("D25").Value = -22
$hist.Range("E25").Value = 1000
$hist.Range("F25").Value = 22000
$hist.Range("G25").Value = "Venta"

$hist.Range("A26").Value = "2025-08-04 20:31:58"
$hist.Range("B26").Value = "Blestro"
$hist.Range("C26").Value = "Nalgas de tom"
$hist.Range("D26").Value = -2
$hist.Range("E26").Value = 100
$hist.Range("F26").Value = 200
$hist.Range("G26").Value = "Venta"

$hist.Range("A27").Value = "2025-08-04 20:32:19"
$hist.Range("B27").Value = "Blestro"
$hist.Range("C27").Value = "Vodka Smirnoff"
$hist.Range("D27").Value = -47
$hist.Range("E27").Value = 5200
$hist.Range("F27").Value = 244400
$hist.Range("G27").Value = "Venta"

$hist.Range("A28").Value = "2025-08-04 20:32:29"
$hist.Range("B28").Value = "Blestro"
$hist.Range("C28").Value = "Fernet Branca"
$hist.Range("D28").Value = -47
$hist.Range("E28").Value = 4500
$hist.Range("F28").Value = 211500
$hist.Range("G28").Value = "Venta"
